# Replace the division problem texts in the table according to the diff.
# Each "old" text is unique within the document, so Find/Replace All is safe.

$d = $word.ActiveDocument

$pairs = @(
    @("53÷6=8, 5", "57÷9=6, 3"),
    @("69÷3=23, 0", "94÷4=23, 2"),
    @("17÷7=2, 3", "88÷9=9, 7"),
    @("45÷4=11, 1", "45÷2=22, 1"),
    @("59÷8=7, 3", "59÷2=29, 1"),
    @("28÷5=5, 3", "36÷6=6, 0"),
    @("17÷6=2, 5", "75÷8=9, 3"),
    @("58÷7=8, 2", "35÷9=3, 8"),
    @("35÷5=7, 0", "69÷9=7, 6"),
    @("42÷7=6, 0", "84÷9=9, 3"),
    @("12÷9=1, 3", "70÷7=10, 0"),
    @("55÷5=11, 0", "50÷2=25, 0"),
    @("23÷4=5, 3", "96÷9=10, 6"),
    @("80÷9=8, 8", "57÷9=6, 3"),
    @("46÷2=23, 0", "51÷9=5, 6"),
    @("21÷3=7, 0", "31÷8=3, 7"),
    @("93÷6=15, 3", "18÷7=2, 4"),
    @("54÷3=18, 0", "82÷7=11, 5"),
    @("49÷9=5, 4", "81÷8=10, 1"),
    @("71÷7=10, 1", "56÷8=7, 0"),
    @("15÷2=7, 1", "10÷7=1, 3"),
    @("60÷8=7, 4", "11÷6=1, 5"),
    @("46÷6=7, 4", "70÷6=11, 4"),
    @("68÷3=22, 2", "75÷9=8, 3"),
    @("20÷6=3, 2", "66÷8=8, 2")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
